$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C3: "Rocky@123" -> "Rocky123" (hyperlinked cell). Update cell text first, then
# restore the original display text onto the existing hyperlink (Excel records
# the old displayed text as the hyperlink's "display" attribute when the cell
# text diverges from it).
$ws.Range("C3").Value = "Rocky123"

# C10 / D10: "Anil123"/"Ajay123" -> "Anil"/"Ajay" (now match C8/D8)
$ws.Range("C10").Value = "Anil"
$ws.Range("D10").Value = "Ajay"

# Restore the hyperlink's display text in place (must go through the pipeline
# so the existing hyperlink entry is updated rather than a new one appended).
$h = $ws.Hyperlinks | Select-Object -First 1
$h.TextToDisplay = "Rocky@123"
